# Edit LOM3108.docx to match the target revision.
# Strategy: the document always has 18 top-level paragraphs, both before and
# after the edit (no paragraphs are added or removed - only the text content
# of certain non-heading paragraphs changes). We therefore address each
# paragraph by its stable 1-based index and rewrite its content explicitly.
#
# NOTE: Paragraph.Range's .Text setter in this host only replaces the first
# run of a multi-run paragraph, so whole-paragraph rewrites instead build an
# explicit Range via $d.Range(start, end), which reliably replaces the full
# span. For the one paragraph that must keep some runs (the bold
# "Método:"/"Critério:"/"Norma de recuperação:" labels) untouched, we use
# targeted Find/Replace across the unwanted text instead.

$d = $word.ActiveDocument
$LB = [char]11   # manual line break -> <w:br/>

# ---- source text fragments (verbatim from the original document) ----
$OBJETIVOS_TEXT = "Promover a formação do engenheiro de materiais sob o ponto de vista do desenvolvimento de competências gerais e específicas. Aplicar e integrar conhecimentos adquiridos às demais disciplinas do curso de Engenharia de Materiais, desenvolvendo competências técnicas relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes) e socioemocionais, num ambiente de aprendizagem colaborativa."
$DOC_CASSIUS    = "3586455 - Cassius Olivio Figueiredo Terra Ruchert"
$DOC_MARIA      = "7459752 - Maria Ismenia Sodero Toledo Faria"
$DOC_SANDRA     = "2166002 - Sandra Giacomin Schneider"
$DOC_SEBASTIAO  = "1922320 - Sebastiao Ribeiro"
$RESUMO_TEXT    = "Introdução a métodos de projeto: metodologias de projeto Design Thinking, Projeto Modelo Canvas e Lean Startup. Ciclo de vida de projeto PDCA Proposta e simulação de pequeno projeto de Engenharia. Definição do problema e formação de alternativas de solução. Estabelecimento de critérios. Escolha e avaliação de soluções. Especificação da solução. Prática de escrita científica."

$PROG1 = "1.Introdução ao projeto em Engenharia: o que é projeto em engenharia e por que projetar? Metodologias de projeto; etapas de elaboração de projeto;"
$PROG2 = "2.Metodologia de projeto focada no ser humano Design Thinking. Entendimento do duplo diamante da inovação. Etapas do Design Thinking: empatia, definição do problema, ideação, prototipação do plano e teste do produto;"
$PROG3 = "3. Processo de melhoria contínua Kaizen. Ciclo de vida de projeto PDCA (Plan-Do-Check-Act): Planejar-Desenvolver-Checar-Agir;"
$PROG4 = "4.Métodos e normas para redação de textos científicos;"
$PROG5 = "5.Desenvolvimento de um projeto temático, compreendendo: definição do problema e formação de alternativas de solução; estabelecimento de critérios; escolha e avaliação de soluções; especificação da solução;"
$PROG6 = "6.Noções de planejamento e gestão de projetos; organização do tempo; técnicas para a realização de apresentações; noções de aprendizagem baseada em projetos; trabalho em grupo, equipes e times"
$PROG7 = "7.Tutoria de projetos"

$METODO1 = "O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros."
$METODO2 = "Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão."
$METODO3 = "Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto."
$METODO4 = "As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas."

$CRITERIO1 = "A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros."
$CRITERIO2 = "O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina."

$NAOHA = "não há"

$BIB1 = "- BAZZO, Walter; PEREIRA, Luiz T.V. Introdução à Engenharia, 3a. edição. Florianópolis: Editora da UFSC, 2013."
$BIB2 = "- COCIAN, Luis Fernando Espinosa. Introdução à Engenharia. Porto Alegre: Bookman, 2017."
$BIB3 = "- BENNETT, Ronald; MILLAM, Elaine. Liderança para engenheiros. Porto Alegre: AMGH, 2014."
$BIB4 = "- ALEXANDER, C. K.; WATSON, J. A. Habilidades para uma carreira de sucesso na engenharia, Porto Alegre: AMGH Editora, 2015."
$BIB5 = "- MCCAHAN, S.; ANDERSON, P.; KORTSCHOT, M.; WEISS, P.; WOODHOUSE, K. Projetos de Engenharia: uma introdução. 1ª edição. -Rio de Janeiro: LTC, 2017."

function Set-ParaText($index, $text) {
    $p = $d.Paragraphs($index)
    $s = $p.Range.Start
    $e = $p.Range.End
    $r = $d.Range($s, $e)
    $r.Text = $text
}

# Paragraph 6: "Objetivos" body text -> the short "Programa resumido" text
Set-ParaText 6 $RESUMO_TEXT

# Paragraph 8: "Docente(s)" bullet list -> objectives text + full program list
#              + the "Método" narrative + the "Critério" narrative
$p8 = $OBJETIVOS_TEXT + $LB + $PROG1 + $LB + $PROG2 + $LB + $PROG3 + $LB + $PROG4 + $LB + $PROG5 + $LB + $PROG6 + $LB + $PROG7 + $LB + $METODO1 + $LB + $METODO2 + $LB + $METODO3 + $LB + $METODO4 + $LB + $CRITERIO1 + $LB + $CRITERIO2
Set-ParaText 8 $p8

# Paragraph 10: "Programa resumido" body text -> "não há"
Set-ParaText 10 $NAOHA

# Paragraph 12: "Programa" body text (numbered list) -> bibliography list
$p12 = $BIB1 + $LB + $BIB2 + $LB + $BIB3 + $LB + $BIB4 + $LB + $BIB5
Set-ParaText 12 $p12

# Paragraph 14: "Avaliação" bullet list - keep the bold "Método:"/"Critério:"/
# "Norma de recuperação:" labels in place, only swap the narrative text after
# each label for the three docente lines. Scope the Find to paragraph 14's
# own range (re-fetched fresh each time since earlier replacements shift
# character offsets) so it cannot match look-alike text that now lives in
# other paragraphs (e.g. paragraph 8 also contains METODO1.. after its
# rewrite above).
$old14a = $METODO1 + $LB + $METODO2 + $LB + $METODO3 + $LB + $METODO4
$p14 = $d.Paragraphs(14)
$rng14 = $d.Range($p14.Range.Start, $p14.Range.End)
$rng14.Find.Execute($old14a, $true, $false, $false, $false, $false, $true, 1, $false, $DOC_CASSIUS, 2) | Out-Null

$old14b = $CRITERIO1 + $LB + $CRITERIO2
$p14 = $d.Paragraphs(14)
$rng14 = $d.Range($p14.Range.Start, $p14.Range.End)
$rng14.Find.Execute($old14b, $true, $false, $false, $false, $false, $true, 1, $false, $DOC_MARIA, 2) | Out-Null

$p14 = $d.Paragraphs(14)
$rng14 = $d.Range($p14.Range.Start, $p14.Range.End)
$rng14.Find.Execute($NAOHA, $true, $false, $false, $false, $false, $true, 1, $false, $DOC_SANDRA, 2) | Out-Null

# Paragraph 16: "Bibliografia" body text -> remaining docente line
Set-ParaText 16 $DOC_SEBASTIAO
